# Lab Exam 03 grading workbook — fill in the "Points for grading" (column E)
# scores for the Customer Class and Product Class rubric blocks, matching
# column D (full marks awarded), and leave the cursor/selection on the
# newly-completed "Customer Class" total cell (E15) instead of the old
# scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Customer Class block (rows 3-6): award full points in column E ---
$ws.Range("E3").Value  = 1
$ws.Range("E4").Value  = 2
$ws.Range("E5").Value  = 2
$ws.Range("E6").Value  = 2

# --- Product Class block (rows 10-14): award full points in column E ---
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Row/Subtotal formulas (E7, E15, E26, E31, E35, E38) already reference
# these cells and recalculate automatically.

# Move the selection to reflect where grading left off, and drop the old
# scrolled-down viewport (previously topLeftCell = A24, selection = C44).
[void]$ws.Range("E15").Select()
